$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1 - match the formatting already used by the
# other header cells (bold, centered, bordered) by copying G1's format.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# New data value for H2 (row 2's "Save" figure).
$ws.Range("H2").Value = 0
